$wb = $excel.ActiveWorkbook

# Sheet: general
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 773.467865015291
$ws.Range("B4").Value = 0.01800012588500977
$ws.Range("B6").Value = 45.77786501528455
$ws.Range("B9").Value = 154.4000000000063
$ws.Range("B10").Value = 573.29

# Sheet: x
$ws = $wb.Worksheets.Item("x")
$ws.Range("B5").Value = 11
$ws.Range("B8").Value = 10
$ws.Range("B11").Value = 13
$ws.Range("B12").Value = 12
$ws.Range("B13").Value = 7
$ws.Range("B14").Value = 1
$ws.Range("B3").Value = 2

# Sheet: U
$ws = $wb.Worksheets.Item("U")
$ws.Range("B8").Value = 3

# Sheet: TBar
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 22.61068956408072
$ws.Range("B4").Value = 10
$ws.Range("B9").Value = 23.22876137241512
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 30
$ws.Range("B14").Value = 25.21630137166173
$ws.Range("B15").Value = 28.21630585843227

# Sheet: y
$ws = $wb.Worksheets.Item("y")
$ws.Range("B2").Value = 11
$ws.Range("B3").Value = 11
$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 11

# Sheet: Q
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 285.6150000000008
$ws.Range("C8").Value = 294.9850000000008
$ws.Range("C9").Value = 290.9700000000008
$ws.Range("C10").Value = 295.4550000000008
$ws.Range("C11").Value = 285.1750000000008
$ws.Range("C12").Value = 67.77500000000072
$ws.Range("C13").Value = 73.77000000000072
$ws.Range("C14").Value = 74.03500000000074
$ws.Range("C15").Value = 73.04500000000073
$ws.Range("C16").Value = 72.66500000000073
$ws.Range("C22").Value = 117.7200000000008
$ws.Range("C23").Value = 118.7150000000007
$ws.Range("C24").Value = 115.9700000000008
$ws.Range("C25").Value = 115.6450000000007
$ws.Range("C26").Value = 116.4550000000007
$ws.Range("C37").Value = 236.3350000000021
$ws.Range("C38").Value = 246.4550000000021
$ws.Range("C39").Value = 231.7250000000021
$ws.Range("C40").Value = 253.5450000000021
$ws.Range("C41").Value = 239.25
$ws.Range("C47").Value = 153.2600000000012
$ws.Range("C48").Value = 161.7350000000012
$ws.Range("C49").Value = 153.75
$ws.Range("C50").Value = 163.7750000000012
$ws.Range("C51").Value = 157.3950000000012
$ws.Range("C52").Value = 175.67
$ws.Range("C53").Value = 177.395
$ws.Range("C54").Value = 177.7
$ws.Range("C55").Value = 176.3
$ws.Range("C56").Value = 169.37
$ws.Range("C57").Value = 285.6150000000008
$ws.Range("C58").Value = 294.9850000000008
$ws.Range("C59").Value = 290.9700000000008
$ws.Range("C60").Value = 295.4550000000008
$ws.Range("C61").Value = 285.1750000000008
$ws.Range("C62").Value = 212.0549999999987
$ws.Range("C63").Value = 215.8299999999987
$ws.Range("C64").Value = 177.0399999999987
$ws.Range("C65").Value = 198
$ws.Range("C66").Value = 184.7
$ws.Range("C67").Value = 236.3350000000021
$ws.Range("C68").Value = 246.4550000000021
$ws.Range("C69").Value = 231.7250000000021
$ws.Range("C70").Value = 253.5450000000021
$ws.Range("C71").Value = 239.25

# Sheet: R
$ws = $wb.Worksheets.Item("R")
$ws.Range("C2").Value = 10.61500000000058
$ws.Range("C3").Value = 19.98500000000065
$ws.Range("C4").Value = 15.97000000000064
$ws.Range("C5").Value = 20.45500000000069
$ws.Range("C6").Value = 10.1750000000006
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0

# Sheet: L
$ws = $wb.Worksheets.Item("L")
$ws.Range("C7").Value = 13.7
$ws.Range("C8").Value = 6.91
$ws.Range("C9").Value = 10.68
$ws.Range("C10").Value = 7.39
$ws.Range("C11").Value = 14.68
$ws.Range("C32").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("C35").Value = 0
$ws.Range("C36").Value = 0
